$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-order match data (columns F:V) within matchday groups that share the same data_partida (col E) ---
# Each group is right-rotated by one row (last row's F:V content moves to the first row of the group).

# Rows 4-6
$row4 = @{}
for ($c = 6; $c -le 22; $c++) { $row4[$c] = $ws.Cells.Item(4, $c).Value() }
$row5 = @{}
for ($c = 6; $c -le 22; $c++) { $row5[$c] = $ws.Cells.Item(5, $c).Value() }
$row6 = @{}
for ($c = 6; $c -le 22; $c++) { $row6[$c] = $ws.Cells.Item(6, $c).Value() }
for ($c = 6; $c -le 22; $c++) { $ws.Cells.Item(4, $c).Value = $row6[$c] }
for ($c = 6; $c -le 22; $c++) { $ws.Cells.Item(5, $c).Value = $row4[$c] }
for ($c = 6; $c -le 22; $c++) { $ws.Cells.Item(6, $c).Value = $row5[$c] }

# Rows 42-43
$row42 = @{}
for ($c = 6; $c -le 22; $c++) { $row42[$c] = $ws.Cells.Item(42, $c).Value() }
$row43 = @{}
for ($c = 6; $c -le 22; $c++) { $row43[$c] = $ws.Cells.Item(43, $c).Value() }
for ($c = 6; $c -le 22; $c++) { $ws.Cells.Item(42, $c).Value = $row43[$c] }
for ($c = 6; $c -le 22; $c++) { $ws.Cells.Item(43, $c).Value = $row42[$c] }

# Rows 56-59
$row56 = @{}
for ($c = 6; $c -le 22; $c++) { $row56[$c] = $ws.Cells.Item(56, $c).Value() }
$row57 = @{}
for ($c = 6; $c -le 22; $c++) { $row57[$c] = $ws.Cells.Item(57, $c).Value() }
$row58 = @{}
for ($c = 6; $c -le 22; $c++) { $row58[$c] = $ws.Cells.Item(58, $c).Value() }
$row59 = @{}
for ($c = 6; $c -le 22; $c++) { $row59[$c] = $ws.Cells.Item(59, $c).Value() }
for ($c = 6; $c -le 22; $c++) { $ws.Cells.Item(56, $c).Value = $row59[$c] }
for ($c = 6; $c -le 22; $c++) { $ws.Cells.Item(57, $c).Value = $row56[$c] }
for ($c = 6; $c -le 22; $c++) { $ws.Cells.Item(58, $c).Value = $row57[$c] }
for ($c = 6; $c -le 22; $c++) { $ws.Cells.Item(59, $c).Value = $row58[$c] }

# Rows 73-76
$row73 = @{}
for ($c = 6; $c -le 22; $c++) { $row73[$c] = $ws.Cells.Item(73, $c).Value() }
$row74 = @{}
for ($c = 6; $c -le 22; $c++) { $row74[$c] = $ws.Cells.Item(74, $c).Value() }
$row75 = @{}
for ($c = 6; $c -le 22; $c++) { $row75[$c] = $ws.Cells.Item(75, $c).Value() }
$row76 = @{}
for ($c = 6; $c -le 22; $c++) { $row76[$c] = $ws.Cells.Item(76, $c).Value() }
for ($c = 6; $c -le 22; $c++) { $ws.Cells.Item(73, $c).Value = $row76[$c] }
for ($c = 6; $c -le 22; $c++) { $ws.Cells.Item(74, $c).Value = $row73[$c] }
for ($c = 6; $c -le 22; $c++) { $ws.Cells.Item(75, $c).Value = $row74[$c] }
for ($c = 6; $c -le 22; $c++) { $ws.Cells.Item(76, $c).Value = $row75[$c] }

# Rows 79-80
$row79 = @{}
for ($c = 6; $c -le 22; $c++) { $row79[$c] = $ws.Cells.Item(79, $c).Value() }
$row80 = @{}
for ($c = 6; $c -le 22; $c++) { $row80[$c] = $ws.Cells.Item(80, $c).Value() }
for ($c = 6; $c -le 22; $c++) { $ws.Cells.Item(79, $c).Value = $row80[$c] }
for ($c = 6; $c -le 22; $c++) { $ws.Cells.Item(80, $c).Value = $row79[$c] }

# --- Append new match rows 92-101 (matchday 11, 10/10/2023 - 29/10/2023 results) ---
$ws.Range("A91:V91").Copy()
$ws.Range("A92:V101").PasteSpecial(-4122)

$ws.Cells.Item(92, 1).Value = 91
$ws.Cells.Item(92, 2).Value = 'england'
$ws.Cells.Item(92, 3).Value = 'premier-league'
$ws.Cells.Item(92, 4).Value = '2023-2024'
$ws.Cells.Item(92, 5).Value = 45226.875
$ws.Cells.Item(92, 6).Value = 'Crystal Palace'
$ws.Cells.Item(92, 7).Value = 1
$ws.Cells.Item(92, 8).Value = 'Tottenham'
$ws.Cells.Item(92, 9).Value = 2
$ws.Cells.Item(92, 10).Value = 3.55
$ws.Cells.Item(92, 11).Value = '10/10/2023 14:02'
$ws.Cells.Item(92, 12).Value = 4.58
$ws.Cells.Item(92, 13).Value = '27/10/2023 20:59'
$ws.Cells.Item(92, 14).Value = 3.69
$ws.Cells.Item(92, 15).Value = '10/10/2023 14:02'
$ws.Cells.Item(92, 16).Value = 3.85
$ws.Cells.Item(92, 17).Value = '27/10/2023 20:50'
$ws.Cells.Item(92, 18).Value = 1.98
$ws.Cells.Item(92, 19).Value = '10/10/2023 14:02'
$ws.Cells.Item(92, 20).Value = 1.81
$ws.Cells.Item(92, 21).Value = '27/10/2023 20:38'
$ws.Cells.Item(92, 22).Value = 'https://www.betexplorer.com/football/england/premier-league/crystal-palace-tottenham/zZOEqy2t/'

$ws.Cells.Item(93, 1).Value = 92
$ws.Cells.Item(93, 2).Value = 'england'
$ws.Cells.Item(93, 3).Value = 'premier-league'
$ws.Cells.Item(93, 4).Value = '2023-2024'
$ws.Cells.Item(93, 5).Value = 45227.5625
$ws.Cells.Item(93, 6).Value = 'Chelsea'
$ws.Cells.Item(93, 7).Value = 0
$ws.Cells.Item(93, 8).Value = 'Brentford'
$ws.Cells.Item(93, 9).Value = 2
$ws.Cells.Item(93, 10).Value = 1.72
$ws.Cells.Item(93, 11).Value = '10/10/2023 14:02'
$ws.Cells.Item(93, 12).Value = 1.65
$ws.Cells.Item(93, 13).Value = '28/10/2023 13:28'
$ws.Cells.Item(93, 14).Value = 3.93
$ws.Cells.Item(93, 15).Value = '10/10/2023 14:02'
$ws.Cells.Item(93, 16).Value = 4.05
$ws.Cells.Item(93, 17).Value = '28/10/2023 13:25'
$ws.Cells.Item(93, 18).Value = 5.01
$ws.Cells.Item(93, 19).Value = '10/10/2023 14:02'
$ws.Cells.Item(93, 20).Value = 5.6
$ws.Cells.Item(93, 21).Value = '28/10/2023 13:24'
$ws.Cells.Item(93, 22).Value = 'https://www.betexplorer.com/football/england/premier-league/chelsea-brentford/YPxmlHgP/'

$ws.Cells.Item(94, 1).Value = 93
$ws.Cells.Item(94, 2).Value = 'england'
$ws.Cells.Item(94, 3).Value = 'premier-league'
$ws.Cells.Item(94, 4).Value = '2023-2024'
$ws.Cells.Item(94, 5).Value = 45227.66666666666
$ws.Cells.Item(94, 6).Value = 'Arsenal'
$ws.Cells.Item(94, 7).Value = 5
$ws.Cells.Item(94, 8).Value = 'Sheffield Utd'
$ws.Cells.Item(94, 9).Value = 0
$ws.Cells.Item(94, 10).Value = 1.17
$ws.Cells.Item(94, 11).Value = '10/10/2023 14:22'
$ws.Cells.Item(94, 12).Value = 1.13
$ws.Cells.Item(94, 13).Value = '28/10/2023 14:59'
$ws.Cells.Item(94, 14).Value = 7.34
$ws.Cells.Item(94, 15).Value = '10/10/2023 14:22'
$ws.Cells.Item(94, 16).Value = 9.5
$ws.Cells.Item(94, 17).Value = '28/10/2023 14:59'
$ws.Cells.Item(94, 18).Value = 12.88
$ws.Cells.Item(94, 19).Value = '10/10/2023 14:22'
$ws.Cells.Item(94, 20).Value = 21
$ws.Cells.Item(94, 21).Value = '28/10/2023 14:59'
$ws.Cells.Item(94, 22).Value = 'https://www.betexplorer.com/football/england/premier-league/arsenal-sheffield-utd/p40dRgnC/'

$ws.Cells.Item(95, 1).Value = 94
$ws.Cells.Item(95, 2).Value = 'england'
$ws.Cells.Item(95, 3).Value = 'premier-league'
$ws.Cells.Item(95, 4).Value = '2023-2024'
$ws.Cells.Item(95, 5).Value = 45227.66666666666
$ws.Cells.Item(95, 6).Value = 'Bournemouth'
$ws.Cells.Item(95, 7).Value = 2
$ws.Cells.Item(95, 8).Value = 'Burnley'
$ws.Cells.Item(95, 9).Value = 1
$ws.Cells.Item(95, 10).Value = 2.22
$ws.Cells.Item(95, 11).Value = '10/10/2023 14:20'
$ws.Cells.Item(95, 12).Value = 2.13
$ws.Cells.Item(95, 13).Value = '28/10/2023 15:57'
$ws.Cells.Item(95, 14).Value = 3.45
$ws.Cells.Item(95, 15).Value = '10/10/2023 14:20'
$ws.Cells.Item(95, 16).Value = 3.65
$ws.Cells.Item(95, 17).Value = '28/10/2023 15:59'
$ws.Cells.Item(95, 18).Value = 3.18
$ws.Cells.Item(95, 19).Value = '10/10/2023 14:20'
$ws.Cells.Item(95, 20).Value = 3.51
$ws.Cells.Item(95, 21).Value = '28/10/2023 15:57'
$ws.Cells.Item(95, 22).Value = 'https://www.betexplorer.com/football/england/premier-league/bournemouth-burnley/W0dhSZW5/'

$ws.Cells.Item(96, 1).Value = 95
$ws.Cells.Item(96, 2).Value = 'england'
$ws.Cells.Item(96, 3).Value = 'premier-league'
$ws.Cells.Item(96, 4).Value = '2023-2024'
$ws.Cells.Item(96, 5).Value = 45227.77083333334
$ws.Cells.Item(96, 6).Value = 'Wolves'
$ws.Cells.Item(96, 7).Value = 2
$ws.Cells.Item(96, 8).Value = 'Newcastle'
$ws.Cells.Item(96, 9).Value = 2
$ws.Cells.Item(96, 10).Value = 4.32
$ws.Cells.Item(96, 11).Value = '10/10/2023 14:02'
$ws.Cells.Item(96, 12).Value = 4.11
$ws.Cells.Item(96, 13).Value = '28/10/2023 18:25'
$ws.Cells.Item(96, 14).Value = 3.84
$ws.Cells.Item(96, 15).Value = '10/10/2023 14:02'
$ws.Cells.Item(96, 16).Value = 3.75
$ws.Cells.Item(96, 17).Value = '28/10/2023 18:28'
$ws.Cells.Item(96, 18).Value = 1.84
$ws.Cells.Item(96, 19).Value = '10/10/2023 14:02'
$ws.Cells.Item(96, 20).Value = 1.93
$ws.Cells.Item(96, 21).Value = '28/10/2023 18:28'
$ws.Cells.Item(96, 22).Value = 'https://www.betexplorer.com/football/england/premier-league/wolves-newcastle-utd/4EGVug15/'

$ws.Cells.Item(97, 1).Value = 96
$ws.Cells.Item(97, 2).Value = 'england'
$ws.Cells.Item(97, 3).Value = 'premier-league'
$ws.Cells.Item(97, 4).Value = '2023-2024'
$ws.Cells.Item(97, 5).Value = 45228.58333333334
$ws.Cells.Item(97, 6).Value = 'West Ham'
$ws.Cells.Item(97, 7).Value = 0
$ws.Cells.Item(97, 8).Value = 'Everton'
$ws.Cells.Item(97, 9).Value = 1
$ws.Cells.Item(97, 10).Value = 1.98
$ws.Cells.Item(97, 11).Value = '10/10/2023 14:02'
$ws.Cells.Item(97, 12).Value = 2.08
$ws.Cells.Item(97, 13).Value = '29/10/2023 13:58'
$ws.Cells.Item(97, 14).Value = 3.79
$ws.Cells.Item(97, 15).Value = '10/10/2023 14:02'
$ws.Cells.Item(97, 16).Value = 3.7
$ws.Cells.Item(97, 17).Value = '29/10/2023 13:58'
$ws.Cells.Item(97, 18).Value = 3.48
$ws.Cells.Item(97, 19).Value = '10/10/2023 14:02'
$ws.Cells.Item(97, 20).Value = 3.62
$ws.Cells.Item(97, 21).Value = '29/10/2023 13:56'
$ws.Cells.Item(97, 22).Value = 'https://www.betexplorer.com/football/england/premier-league/west-ham-everton/QNHRtZoa/'

$ws.Cells.Item(98, 1).Value = 97
$ws.Cells.Item(98, 2).Value = 'england'
$ws.Cells.Item(98, 3).Value = 'premier-league'
$ws.Cells.Item(98, 4).Value = '2023-2024'
$ws.Cells.Item(98, 5).Value = 45228.625
$ws.Cells.Item(98, 6).Value = 'Brighton'
$ws.Cells.Item(98, 7).Value = 1
$ws.Cells.Item(98, 8).Value = 'Fulham'
$ws.Cells.Item(98, 9).Value = 1
$ws.Cells.Item(98, 10).Value = 1.51
$ws.Cells.Item(98, 11).Value = '10/10/2023 14:02'
$ws.Cells.Item(98, 12).Value = 1.64
$ws.Cells.Item(98, 13).Value = '29/10/2023 14:58'
$ws.Cells.Item(98, 14).Value = 4.73
$ws.Cells.Item(98, 15).Value = '10/10/2023 14:02'
$ws.Cells.Item(98, 16).Value = 4.32
$ws.Cells.Item(98, 17).Value = '29/10/2023 14:59'
$ws.Cells.Item(98, 18).Value = 6.2
$ws.Cells.Item(98, 19).Value = '10/10/2023 14:02'
$ws.Cells.Item(98, 20).Value = 5.29
$ws.Cells.Item(98, 21).Value = '29/10/2023 14:59'
$ws.Cells.Item(98, 22).Value = 'https://www.betexplorer.com/football/england/premier-league/brighton-fulham/6Jl5PXGO/'

$ws.Cells.Item(99, 1).Value = 98
$ws.Cells.Item(99, 2).Value = 'england'
$ws.Cells.Item(99, 3).Value = 'premier-league'
$ws.Cells.Item(99, 4).Value = '2023-2024'
$ws.Cells.Item(99, 5).Value = 45228.625
$ws.Cells.Item(99, 6).Value = 'Liverpool'
$ws.Cells.Item(99, 7).Value = 3
$ws.Cells.Item(99, 8).Value = 'Nottingham'
$ws.Cells.Item(99, 9).Value = 0
$ws.Cells.Item(99, 10).Value = 1.23
$ws.Cells.Item(99, 11).Value = '10/10/2023 14:02'
$ws.Cells.Item(99, 12).Value = 1.22
$ws.Cells.Item(99, 13).Value = '29/10/2023 14:51'
$ws.Cells.Item(99, 14).Value = 6.65
$ws.Cells.Item(99, 15).Value = '10/10/2023 14:02'
$ws.Cells.Item(99, 16).Value = 7.31
$ws.Cells.Item(99, 17).Value = '29/10/2023 14:59'
$ws.Cells.Item(99, 18).Value = 9.43
$ws.Cells.Item(99, 19).Value = '10/10/2023 14:02'
$ws.Cells.Item(99, 20).Value = 12.08
$ws.Cells.Item(99, 21).Value = '29/10/2023 14:59'
$ws.Cells.Item(99, 22).Value = 'https://www.betexplorer.com/football/england/premier-league/liverpool-nottingham/IcEJreHn/'

$ws.Cells.Item(100, 1).Value = 99
$ws.Cells.Item(100, 2).Value = 'england'
$ws.Cells.Item(100, 3).Value = 'premier-league'
$ws.Cells.Item(100, 4).Value = '2023-2024'
$ws.Cells.Item(100, 5).Value = 45228.625
$ws.Cells.Item(100, 6).Value = 'Aston Villa'
$ws.Cells.Item(100, 7).Value = 3
$ws.Cells.Item(100, 8).Value = 'Luton'
$ws.Cells.Item(100, 9).Value = 1
$ws.Cells.Item(100, 10).Value = 1.38
$ws.Cells.Item(100, 11).Value = '10/10/2023 14:32'
$ws.Cells.Item(100, 12).Value = 1.34
$ws.Cells.Item(100, 13).Value = '29/10/2023 14:55'
$ws.Cells.Item(100, 14).Value = 4.92
$ws.Cells.Item(100, 15).Value = '10/10/2023 14:32'
$ws.Cells.Item(100, 16).Value = 5.85
$ws.Cells.Item(100, 17).Value = '29/10/2023 14:55'
$ws.Cells.Item(100, 18).Value = 7.42
$ws.Cells.Item(100, 19).Value = '10/10/2023 14:32'
$ws.Cells.Item(100, 20).Value = 8.55
$ws.Cells.Item(100, 21).Value = '29/10/2023 14:57'
$ws.Cells.Item(100, 22).Value = 'https://www.betexplorer.com/football/england/premier-league/aston-villa-luton/SSk1QD1I/'

$ws.Cells.Item(101, 1).Value = 100
$ws.Cells.Item(101, 2).Value = 'england'
$ws.Cells.Item(101, 3).Value = 'premier-league'
$ws.Cells.Item(101, 4).Value = '2023-2024'
$ws.Cells.Item(101, 5).Value = 45228.6875
$ws.Cells.Item(101, 6).Value = 'Manchester Utd'
$ws.Cells.Item(101, 7).Value = 0
$ws.Cells.Item(101, 8).Value = 'Manchester City'
$ws.Cells.Item(101, 9).Value = 3
$ws.Cells.Item(101, 10).Value = 3.42
$ws.Cells.Item(101, 11).Value = '10/10/2023 14:02'
$ws.Cells.Item(101, 12).Value = 5.6
$ws.Cells.Item(101, 13).Value = '29/10/2023 16:17'
$ws.Cells.Item(101, 14).Value = 3.72
$ws.Cells.Item(101, 15).Value = '10/10/2023 14:02'
$ws.Cells.Item(101, 16).Value = 4.3
$ws.Cells.Item(101, 17).Value = '29/10/2023 16:17'
$ws.Cells.Item(101, 18).Value = 2.02
$ws.Cells.Item(101, 19).Value = '10/10/2023 14:02'
$ws.Cells.Item(101, 20).Value = 1.61
$ws.Cells.Item(101, 21).Value = '29/10/2023 16:17'
$ws.Cells.Item(101, 22).Value = 'https://www.betexplorer.com/football/england/premier-league/manchester-united-manchester-city/W0INsFWh/'

Write-Output "edit applied"
